# Scheduled market-data refresh: update Leve profit computation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1480.4445
$ws.Range("I43").Value = 1164.9
$ws.Range("J43").Value = 1874.875
$ws.Range("K43").Value = 1164.9
$ws.Range("L43").Value = 1874.875
$ws.Range("M43").Value = -1095.9
$ws.Range("N43").Value = -2012.875

# Row 51
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 1300
$ws.Range("K51").Value = 1300
$ws.Range("M51").Value = -816

# Row 106
$ws.Range("H106").Value = 28574996
$ws.Range("I106").Value = 33337050
$ws.Range("J106").Value = 2666.6667
$ws.Range("K106").Value = 33337050
$ws.Range("L106").Value = 2666.6667
$ws.Range("M106").Value = -33336419
$ws.Range("N106").Value = -3928.6667

# Row 138
$ws.Range("H138").Value = 1306.17
$ws.Range("I138").Value = 597.0833
$ws.Range("J138").Value = 1960.7115
$ws.Range("K138").Value = 1791.2499
$ws.Range("L138").Value = 5882.1345
$ws.Range("M138").Value = 3348.7501
$ws.Range("N138").Value = -16162.1345

# Row 141
$ws.Range("H141").Value = 2842.1428
$ws.Range("I141").Value = 952.3611
$ws.Range("J141").Value = 14180.833
$ws.Range("K141").Value = 2857.0833
$ws.Range("L141").Value = 42542.499
$ws.Range("M141").Value = 2322.9167
$ws.Range("N141").Value = -52902.499

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5377893.5
$ws.Range("I61").Value = 6290509.5
$ws.Range("J61").Value = 3600.2222
$ws.Range("K61").Value = 6290509.5
$ws.Range("L61").Value = 3600.2222
$ws.Range("M61").Value = -6290297.5
$ws.Range("N61").Value = -4024.2222

# Row 74
$ws.Range("H74").Value = 1502.6792
$ws.Range("I74").Value = 839.8148
$ws.Range("J74").Value = 2191.0386
$ws.Range("K74").Value = 839.8148
$ws.Range("L74").Value = 2191.0386
$ws.Range("M74").Value = 34.18520000000001
$ws.Range("N74").Value = -3939.0386

# Row 77
$ws.Range("H77").Value = 1502.6792
$ws.Range("I77").Value = 839.8148
$ws.Range("J77").Value = 2191.0386
$ws.Range("K77").Value = 4199.074
$ws.Range("L77").Value = 10955.193
$ws.Range("M77").Value = 168.9260000000004
$ws.Range("N77").Value = -19691.193

# Row 122
$ws.Range("H122").Value = 44358.26
$ws.Range("J122").Value = 1657
$ws.Range("L122").Value = 4971
$ws.Range("N122").Value = -9871

# Row 132
$ws.Range("H132").Value = 2166.2837
$ws.Range("I132").Value = 1331.431
$ws.Range("J132").Value = 5192.625
$ws.Range("K132").Value = 3994.293
$ws.Range("L132").Value = 15577.875
$ws.Range("M132").Value = -1464.293
$ws.Range("N132").Value = -20637.875

# Row 136
$ws.Range("H136").Value = 5377893.5
$ws.Range("I136").Value = 6290509.5
$ws.Range("J136").Value = 3600.2222
$ws.Range("K136").Value = 18871528.5
$ws.Range("L136").Value = 10800.6666
$ws.Range("M136").Value = -18868978.5
$ws.Range("N136").Value = -15900.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1263.7693
$ws.Range("I16").Value = 1011.7273
$ws.Range("K16").Value = 1011.7273
$ws.Range("M16").Value = -724.7273

# Row 31
$ws.Range("H31").Value = 4502.852
$ws.Range("I31").Value = 1406.6072
$ws.Range("J31").Value = 7837.269
$ws.Range("K31").Value = 1406.6072
$ws.Range("L31").Value = 7837.269
$ws.Range("M31").Value = -1111.6072
$ws.Range("N31").Value = -8427.269

# Row 34
$ws.Range("H34").Value = 4502.852
$ws.Range("I34").Value = 1406.6072
$ws.Range("J34").Value = 7837.269
$ws.Range("K34").Value = 1406.6072
$ws.Range("L34").Value = 7837.269
$ws.Range("M34").Value = -1204.6072
$ws.Range("N34").Value = -8241.269

# Row 58
$ws.Range("H58").Value = 959.36957
$ws.Range("I58").Value = 660.3333
$ws.Range("J58").Value = 1718.4615
$ws.Range("K58").Value = 660.3333
$ws.Range("L58").Value = 1718.4615
$ws.Range("M58").Value = -457.3333
$ws.Range("N58").Value = -2124.4615

# Row 113
$ws.Range("H113").Value = 1263.7693
$ws.Range("I113").Value = 1011.7273
$ws.Range("K113").Value = 1011.7273
$ws.Range("M113").Value = 1158.2727

# Row 122
$ws.Range("H122").Value = 1708.3055
$ws.Range("I122").Value = 1368.75
$ws.Range("K122").Value = 4106.25
$ws.Range("M122").Value = -1656.25

# Row 132
$ws.Range("H132").Value = 3969677
$ws.Range("I132").Value = 1218.3103
$ws.Range("J132").Value = 12822392
$ws.Range("K132").Value = 3654.9309
$ws.Range("L132").Value = 38467176
$ws.Range("M132").Value = -1124.9309
$ws.Range("N132").Value = -38472236

# Row 134
$ws.Range("H134").Value = 2892.0635
$ws.Range("I134").Value = 2948
$ws.Range("J134").Value = 2696.2856
$ws.Range("K134").Value = 8844
$ws.Range("L134").Value = 8088.8568
$ws.Range("M134").Value = -6309
$ws.Range("N134").Value = -13158.8568

# Row 136
$ws.Range("H136").Value = 959.36957
$ws.Range("I136").Value = 660.3333
$ws.Range("J136").Value = 1718.4615
$ws.Range("K136").Value = 1980.9999
$ws.Range("L136").Value = 5155.3845
$ws.Range("M136").Value = 569.0001
$ws.Range("N136").Value = -10255.3845

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1350.6666
$ws.Range("I5").Value = 383.7143
$ws.Range("J5").Value = 2063.158
$ws.Range("K5").Value = 1151.1429
$ws.Range("L5").Value = 6189.474
$ws.Range("M5").Value = -1039.1429
$ws.Range("N5").Value = -6413.474

# Row 122
$ws.Range("H122").Value = 3154.9
$ws.Range("I122").Value = 436.73914
$ws.Range("J122").Value = 6832.4116
$ws.Range("K122").Value = 3930.65226
$ws.Range("L122").Value = 61491.7044
$ws.Range("M122").Value = -1480.65226
$ws.Range("N122").Value = -66391.7044

# Row 131
$ws.Range("H131").Value = 3040.611
$ws.Range("I131").Value = 358.83334
$ws.Range("J131").Value = 4381.5
$ws.Range("K131").Value = 1076.50002
$ws.Range("L131").Value = 13144.5
$ws.Range("M131").Value = 3963.49998
$ws.Range("N131").Value = -23224.5

# Row 135
$ws.Range("H135").Value = 1350.6666
$ws.Range("I135").Value = 383.7143
$ws.Range("J135").Value = 2063.158
$ws.Range("K135").Value = 3453.4287
$ws.Range("L135").Value = 18568.422
$ws.Range("M135").Value = -918.4286999999999
$ws.Range("N135").Value = -23638.422

# Row 139
$ws.Range("H139").Value = 2263.804
$ws.Range("I139").Value = 1155.7142
$ws.Range("J139").Value = 3612.7827
$ws.Range("K139").Value = 3467.1426
$ws.Range("L139").Value = 10838.3481
$ws.Range("M139").Value = 1672.8574
$ws.Range("N139").Value = -21118.3481

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 6101.875
$ws.Range("J29").Value = 6901.143
$ws.Range("L29").Value = 6901.143
$ws.Range("N29").Value = -7481.143

# Row 121
$ws.Range("H121").Value = 69500
$ws.Range("J121").Value = 69500
$ws.Range("L121").Value = 69500
$ws.Range("N121").Value = -72994

# Row 132
$ws.Range("H132").Value = 1807.9419
$ws.Range("I132").Value = 1503.3732
$ws.Range("J132").Value = 2881.9473
$ws.Range("K132").Value = 4510.1196
$ws.Range("L132").Value = 8645.841899999999
$ws.Range("M132").Value = -1980.1196
$ws.Range("N132").Value = -13705.8419

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2828.9546
$ws.Range("I61").Value = 2489.5625
$ws.Range("K61").Value = 2489.5625
$ws.Range("M61").Value = -2287.5625

# Row 68
$ws.Range("H68").Value = 1454.48
$ws.Range("I68").Value = 1418.3112
$ws.Range("J68").Value = 1780
$ws.Range("K68").Value = 1418.3112
$ws.Range("L68").Value = 1780
$ws.Range("M68").Value = -669.3112000000001
$ws.Range("N68").Value = -3278

# Row 69
$ws.Range("H69").Value = 31000
$ws.Range("J69").Value = 31000
$ws.Range("L69").Value = 31000
$ws.Range("N69").Value = -32622

# Row 71
$ws.Range("H71").Value = 1454.48
$ws.Range("I71").Value = 1418.3112
$ws.Range("J71").Value = 1780
$ws.Range("K71").Value = 7091.556
$ws.Range("L71").Value = 8900
$ws.Range("M71").Value = -3347.556
$ws.Range("N71").Value = -16388

# Row 72
$ws.Range("H72").Value = 31000
$ws.Range("J72").Value = 31000
$ws.Range("L72").Value = 93000
$ws.Range("N72").Value = -101112

# Row 113
$ws.Range("H113").Value = 2828.9546
$ws.Range("I113").Value = 2489.5625
$ws.Range("K113").Value = 2489.5625
$ws.Range("M113").Value = -319.5625

# Row 132
$ws.Range("H132").Value = 2586.3235
$ws.Range("I132").Value = 2279.6182
$ws.Range("K132").Value = 6838.8546
$ws.Range("M132").Value = -4308.8546

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 41100.625
$ws.Range("I62").Value = 2700.6667
$ws.Range("J62").Value = 64140.6
$ws.Range("K62").Value = 2700.6667
$ws.Range("L62").Value = 64140.6
$ws.Range("M62").Value = -2076.6667
$ws.Range("N62").Value = -65388.6

# Row 65
$ws.Range("H65").Value = 41100.625
$ws.Range("I65").Value = 2700.6667
$ws.Range("J65").Value = 64140.6
$ws.Range("K65").Value = 13503.3335
$ws.Range("L65").Value = 320703
$ws.Range("M65").Value = -10383.3335
$ws.Range("N65").Value = -326943

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 136
$ws.Range("H136").Value = 1387.52
$ws.Range("I136").Value = 1336.117
$ws.Range("J136").Value = 1559.6086
$ws.Range("K136").Value = 4008.351
$ws.Range("L136").Value = 4678.825800000001
$ws.Range("M136").Value = -1458.351
$ws.Range("N136").Value = -11809.5558
